# Updates the cryptos list (Price / Volume(1h) columns) on sheet1 to the
# latest scraped values, per the Sat May  4 13:39:23 UTC 2024 GitHub
# Actions commit.
#
# Each value is assigned with a leading apostrophe, the same trick used
# when typing a number into the Excel UI and wanting it kept as literal
# text, so Excel stores it as text instead of auto-coercing look-alike
# numbers (e.g. 589.60 -> 589.6, 0.0000255 -> 2.55E-05) -- matching how
# these Price/Volume cells were authored originally (plain text, General
# number format).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.627.78"
$ws.Range("E2").Value = "'  +3.41%  "
$ws.Range("D3").Value = "'3.134.29"
$ws.Range("E3").Value = "'  +2.56%  "
$ws.Range("E4").Value = "'  +0.12%  "
$ws.Range("D5").Value = "'589.60"
$ws.Range("E5").Value = "'  +2.21%  "
$ws.Range("D6").Value = "'146.26"
$ws.Range("E6").Value = "'  +2.42%  "
$ws.Range("E7").Value = "'  -0.01%  "
$ws.Range("D8").Value = "'3.125.71"
$ws.Range("E8").Value = "'  +2.45%  "
$ws.Range("D9").Value = "'0.536"
$ws.Range("E9").Value = "'  +2.14%  "
$ws.Range("D10").Value = "'0.163"
$ws.Range("E10").Value = "'  +17.74%  "
$ws.Range("D11").Value = "'5.69"
$ws.Range("E11").Value = "'  +4.46%  "
$ws.Range("D12").Value = "'0.467"
$ws.Range("E12").Value = "'  +0.89%  "
$ws.Range("D13").Value = "'0.0000255"
$ws.Range("E13").Value = "'  +7.48%  "
$ws.Range("D14").Value = "'35.72"
$ws.Range("E14").Value = "'  +3.12%  "
$ws.Range("E15").Value = "'  +0.06%  "
$ws.Range("D16").Value = "'3.662.81"
$ws.Range("E16").Value = "'  +2.90%  "
$ws.Range("D17").Value = "'63.598.46"
$ws.Range("E17").Value = "'  +3.48%  "
$ws.Range("D18").Value = "'7.15"
$ws.Range("E18").Value = "'  -0.73%  "
$ws.Range("D19").Value = "'3.136.98"
$ws.Range("E19").Value = "'  +2.90%  "
$ws.Range("D20").Value = "'465.97"
$ws.Range("E20").Value = "'  +4.21%  "
$ws.Range("D21").Value = "'14.19"
$ws.Range("E21").Value = "'  +2.44%  "
$ws.Range("D22").Value = "'0.731"
$ws.Range("E22").Value = "'  +0.05%  "
$ws.Range("D23").Value = "'7.49"
$ws.Range("E23").Value = "'  +3.27%  "
$ws.Range("D24").Value = "'13.30"
$ws.Range("E24").Value = "'  -2.24%  "
$ws.Range("D25").Value = "'82.13"
$ws.Range("E25").Value = "'  +0.66%  "
$ws.Range("D27").Value = "'8.58"
$ws.Range("E27").Value = "'  +6.39%  "
$ws.Range("D28").Value = "'2.70"
$ws.Range("E28").Value = "'  +3.04%  "
$ws.Range("E29").Value = "'  +0.24%  "
$ws.Range("E30").Value = "'  -3.63%  "
$ws.Range("D31").Value = "'6.81"
$ws.Range("E31").Value = "'  +5.94%  "
$ws.Range("D32").Value = "'26.93"
$ws.Range("E32").Value = "'  +1.75%  "
$ws.Range("E33").Value = "'  +1.64%  "
$ws.Range("D34").Value = "'0.0₃0861"
$ws.Range("E34").Value = "'  +6.68%  "
$ws.Range("D35").Value = "'2.40"
$ws.Range("E35").Value = "'  +10.38%  "
$ws.Range("E36").Value = "'  +2.60%  "
$ws.Range("D37").Value = "'3.38"
$ws.Range("E37").Value = "'  +14.66%  "
$ws.Range("D38").Value = "'6.13"
$ws.Range("E38").Value = "'  +1.34%  "
$ws.Range("D39").Value = "'50.87"
$ws.Range("E39").Value = "'  +1.74%  "
$ws.Range("D40").Value = "'453.30"
$ws.Range("E40").Value = "'  +10.26%  "
$ws.Range("D41").Value = "'8.69"
$ws.Range("E41").Value = "'  -1.02%  "
$ws.Range("D42").Value = "'0.0374"
$ws.Range("E42").Value = "'  +3.20%  "
$ws.Range("D43").Value = "'2.918.84"
$ws.Range("E43").Value = "'  +4.70%  "
$ws.Range("D44").Value = "'0.279"
$ws.Range("E44").Value = "'  +6.37%  "
$ws.Range("E45").Value = "'  +3.79%  "
$ws.Range("D46").Value = "'2.15"
$ws.Range("E46").Value = "'  +2.97%  "
$ws.Range("D47").Value = "'126.41"
$ws.Range("E47").Value = "'  +2.14%  "
$ws.Range("D49").Value = "'0.111"
$ws.Range("E49").Value = "'  +0.65%  "
$ws.Range("D50").Value = "'24.73"
$ws.Range("E50").Value = "'  +2.92%  "
$ws.Range("D51").Value = "'33.79"
$ws.Range("E51").Value = "'  -10.04%  "
